# Generate Report for Handback
# Swap the "a.md" / "b.md" rows on every sheet: row 2 now reports on b.md
# (still "not in sync"), row 3 now reports on a.md with a fresh handback
# that is "in sync" with en-US.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "b.md"
$ws.Range("B2").Value = "e2e\b.md"
$ws.Range("E2").Value = "Handed back: not in sync with en-US"
$ws.Range("F2").Value = "Handed back: not in sync with en-US"

$ws.Range("A3").Value = "a.md"
$ws.Range("B3").Value = "e2e\a.md"
$ws.Range("E3").Value = "Handed back: in sync with en-US"
$ws.Range("F3").Value = "Handed back: in sync with en-US"
$ws.Range("G3").Value = "2017-02-09 15:11:21"

$i = 0
foreach ($hl in $ws.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) { $hl.TextToDisplay = "e2e\b.md" }
    if ($i -eq 2) { $hl.TextToDisplay = "e2e\a.md" }
}

$ws.Columns.Item(5).ColumnWidth = 32.67
$ws.Columns.Item(6).ColumnWidth = 32.67

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "b.md"
$ws.Range("C2").Value = "Handed back: not in sync with en-US"
$ws.Range("J2").Value = "b.md"

$ws.Range("A3").Value = "a.md"
$ws.Range("H3").Value = "2017-02-09 15:11:02"
$ws.Range("J3").Value = "a.md"
$ws.Range("L3").Value = "2017-02-09 15:12:48"
$ws.Range("M3").Value = "TestHandback_201702091112"

$i = 0
foreach ($hl in $ws.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) { $hl.TextToDisplay = "b.md" }
    if ($i -eq 2) { $hl.TextToDisplay = "b.md" }
    if ($i -eq 3) { $hl.TextToDisplay = "a.md" }
    if ($i -eq 4) { $hl.TextToDisplay = "a.md" }
}

$ws.Columns.Item(3).ColumnWidth = 32.67

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "b.md"
$ws.Range("C2").Value = "Handed back: not in sync with en-US"
$ws.Range("J2").Value = "b.md"

$ws.Range("A3").Value = "a.md"
$ws.Range("H3").Value = "2017-02-09 15:11:21"
$ws.Range("J3").Value = "a.md"
$ws.Range("L3").Value = "2017-02-09 15:13:13"
$ws.Range("M3").Value = "TestHandback_201702091112"

$i = 0
foreach ($hl in $ws.Hyperlinks) {
    $i = $i + 1
    if ($i -eq 1) { $hl.TextToDisplay = "b.md" }
    if ($i -eq 2) { $hl.TextToDisplay = "b.md" }
    if ($i -eq 3) { $hl.TextToDisplay = "a.md" }
    if ($i -eq 4) { $hl.TextToDisplay = "a.md" }
}

$ws.Columns.Item(3).ColumnWidth = 32.67
